$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2849.183
$ws.Range("I15").Value = 2849.183
$ws.Range("K15").Value = 8547.548999999999
$ws.Range("M15").Value = -8378.548999999999

# Row 17
$ws.Range("H17").Value = 945.87805
$ws.Range("J17").Value = 945.87805
$ws.Range("L17").Value = 2837.63415
$ws.Range("N17").Value = -3173.63415

# Row 40
$ws.Range("H40").Value = 983.5
$ws.Range("I40").Value = 901
$ws.Range("K40").Value = 901
$ws.Range("M40").Value = -726

# Row 62
$ws.Range("H62").Value = 2023.1333
$ws.Range("I62").Value = 1614.125
$ws.Range("J62").Value = 2490.5715
$ws.Range("K62").Value = 1614.125
$ws.Range("L62").Value = 2490.5715
$ws.Range("M62").Value = -990.125
$ws.Range("N62").Value = -3738.5715

# Row 65
$ws.Range("H65").Value = 2023.1333
$ws.Range("I65").Value = 1614.125
$ws.Range("J65").Value = 2490.5715
$ws.Range("K65").Value = 8070.625
$ws.Range("L65").Value = 12452.8575
$ws.Range("M65").Value = -4950.625
$ws.Range("N65").Value = -18692.8575

# Row 125
$ws.Range("H125").Value = 1017.7143
$ws.Range("I125").Value = 406.4
$ws.Range("J125").Value = 1208.75
$ws.Range("K125").Value = 3657.6
$ws.Range("L125").Value = 10878.75
$ws.Range("M125").Value = -1197.6
$ws.Range("N125").Value = -15798.75

# Row 132
$ws.Range("H132").Value = 1198274.9
$ws.Range("I132").Value = 3780.75
$ws.Range("J132").Value = 2884619.5
$ws.Range("K132").Value = 11342.25
$ws.Range("L132").Value = 8653858.5
$ws.Range("M132").Value = -8812.25
$ws.Range("N132").Value = -8658918.5

# Row 138
$ws.Range("H138").Value = 2059589.8
$ws.Range("I138").Value = 1427.3182
$ws.Range("K138").Value = 4281.9546
$ws.Range("M138").Value = 858.0454

$ws = $wb.Worksheets.Item("ARM")
# Row 49
$ws.Range("H49").Value = 14999.667
$ws.Range("J49").Value = 14999.667
$ws.Range("L49").Value = 14999.667
$ws.Range("N49").Value = -15519.667

# Row 61
$ws.Range("H61").Value = 45547756
$ws.Range("I61").Value = 62564264
$ws.Range("J61").Value = 170402.33
$ws.Range("K61").Value = 62564264
$ws.Range("L61").Value = 170402.33
$ws.Range("M61").Value = -62564052
$ws.Range("N61").Value = -170826.33

# Row 122
$ws.Range("H122").Value = 2417216.5
$ws.Range("I122").Value = 1777.8918
$ws.Range("K122").Value = 5333.6754
$ws.Range("M122").Value = -2883.6754

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

# Row 132
$ws.Range("H132").Value = 73379.03
$ws.Range("I132").Value = 50907.9
$ws.Range("K132").Value = 152723.7
$ws.Range("M132").Value = -150193.7

# Row 136
$ws.Range("H136").Value = 45547756
$ws.Range("I136").Value = 62564264
$ws.Range("J136").Value = 170402.33
$ws.Range("K136").Value = 187692792
$ws.Range("L136").Value = 511206.99
$ws.Range("M136").Value = -187690242
$ws.Range("N136").Value = -516306.99

$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("N2").Value = 0

# Row 94
$ws.Range("H94").Value = 671.8570999999999
$ws.Range("I94").Value = 492.76923
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 492.76923
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -41.76922999999999
$ws.Range("N94").Value = -3902

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 11000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572

# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 81
$ws.Range("H81").Value = 31659.334
$ws.Range("J81").Value = 31659.334
$ws.Range("L81").Value = 31659.334
$ws.Range("N81").Value = -33655.334

# Row 84
$ws.Range("H84").Value = 31659.334
$ws.Range("J84").Value = 31659.334
$ws.Range("L84").Value = 94978.00199999999
$ws.Range("N84").Value = -104962.002

# Row 112
$ws.Range("H112").Value = 44980
$ws.Range("J112").Value = 44980
$ws.Range("L112").Value = 44980
$ws.Range("N112").Value = -47934

# Row 122
$ws.Range("H122").Value = 1955.2667
$ws.Range("I122").Value = 1539.375
$ws.Range("J122").Value = 2430.5715
$ws.Range("K122").Value = 4618.125
$ws.Range("L122").Value = 7291.7145
$ws.Range("M122").Value = -2168.125
$ws.Range("N122").Value = -12191.7145

# Row 132
$ws.Range("H132").Value = 21453.568
$ws.Range("I132").Value = 1663.8857
$ws.Range("J132").Value = 64743.5
$ws.Range("K132").Value = 4991.6571
$ws.Range("L132").Value = 194230.5
$ws.Range("M132").Value = -2461.6571
$ws.Range("N132").Value = -199290.5

# Row 134
$ws.Range("H134").Value = 19991.586
$ws.Range("I134").Value = 1387.475
$ws.Range("J134").Value = 61334.055
$ws.Range("K134").Value = 4162.424999999999
$ws.Range("L134").Value = 184002.165
$ws.Range("M134").Value = -1627.424999999999
$ws.Range("N134").Value = -189072.165

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 3500
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = -1331
$ws.Range("N36").Value = -12338

# Row 131
$ws.Range("H131").Value = 1052.6865
$ws.Range("I131").Value = 416.25
$ws.Range("J131").Value = 1138.983
$ws.Range("K131").Value = 1248.75
$ws.Range("L131").Value = 3416.949
$ws.Range("M131").Value = 3791.25
$ws.Range("N131").Value = -13496.949

$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Range("H48").Value = 3027
$ws.Range("I48").Value = 3027
$ws.Range("K48").Value = 3027
$ws.Range("M48").Value = -2542

# Row 112
$ws.Range("H112").Value = 22293
$ws.Range("J112").Value = 22293
$ws.Range("L112").Value = 22293
$ws.Range("N112").Value = -24509

# Row 132
$ws.Range("H132").Value = 45086.805
$ws.Range("I132").Value = 34725.566
$ws.Range("J132").Value = 64514.125
$ws.Range("K132").Value = 104176.698
$ws.Range("L132").Value = 193542.375
$ws.Range("M132").Value = -101646.698
$ws.Range("N132").Value = -198602.375

# Row 136
$ws.Range("H136").Value = 11134.866
$ws.Range("J136").Value = 11134.866
$ws.Range("L136").Value = 33404.598
$ws.Range("N136").Value = -38504.598

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 694.1539
$ws.Range("I22").Value = 443.2143
$ws.Range("J22").Value = 834.6799999999999
$ws.Range("K22").Value = 443.2143
$ws.Range("L22").Value = 834.6799999999999
$ws.Range("M22").Value = -148.2143
$ws.Range("N22").Value = -1424.68

# Row 27
$ws.Range("H27").Value = 694.1539
$ws.Range("I27").Value = 443.2143
$ws.Range("J27").Value = 834.6799999999999
$ws.Range("K27").Value = 443.2143
$ws.Range("L27").Value = 834.6799999999999
$ws.Range("M27").Value = -336.2143
$ws.Range("N27").Value = -1048.68

# Row 47
$ws.Range("H47").Value = 9997
$ws.Range("J47").Value = 9997
$ws.Range("L47").Value = 9997
$ws.Range("N47").Value = -10977

# Row 52
$ws.Range("H52").Value = 9997
$ws.Range("J52").Value = 9997
$ws.Range("L52").Value = 9997
$ws.Range("N52").Value = -10463
